$d = $word.ActiveDocument

$replacements = @(
    @("2024-05-10 Friday", "2024-05-11 Saturday"),
    @("84×37=", "51×40="),
    @("61×74=", "48×90="),
    @("49×26=", "17×34="),
    @("71×88=", "81×66="),
    @("54×26=", "97×42="),
    @("97×23=", "97×93="),
    @("18×77=", "85×72="),
    @("82×99=", "90×51="),
    @("22×68=", "58×41="),
    @("93×30=", "39×70="),
    @("99×66=", "86×91="),
    @("12×85=", "74×40="),
    @("77×78=", "38×57="),
    @("27×59=", "18×42="),
    @("64×16=", "54×57="),
    @("73×75=", "47×78="),
    @("22×45=", "21×73="),
    @("38×77=", "53×39="),
    @("76×35=", "54×94="),
    @("74×47=", "62×40="),
    @("21×34=", "41×71="),
    @("37×55=", "42×50="),
    @("31×69=", "44×25="),
    @("79×93=", "43×63="),
    @("26×19=", "86×20=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
